$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "STATION"
$ws.Range("B2").Value = "NAME"
$ws.Range("C2").Value = "NETID"
$ws.Range("D2").Value = "LAT"
$ws.Range("E2").Value = "LON"

$ws.Range("F2").Select()
